# Daily attendance processing - 2025-11-16 22:20:09
# Normalize the "Recorded By" (column G) entries so that the "System" /
# "system" marker is listed before the user/email identifier(s), instead
# of after them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact before -> after mappings observed for the "Recorded By" column.
$map = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "system, backup@backdoor.com, System" = "system, System, backup@backdoor.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
